# Apply the attendance-sheet edits described by the diff.
#
# Workbook has a single visible sheet named "August" (internal codeName
# Sheet2). Columns D/E/F/G/H hold attendance marks for the class dates in
# row 10 (3rd/10th/17th/24th); column I totals a student's row, column J
# turns that total into a percentage of the "Total Classes" row (row 12).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Scroll position / selection on reopen -------------------------------
# Excel records the window's top-left visible cell and the last selection
# for the sheet. The selection is reachable through Range.Select(); make a
# best-effort attempt at the scroll position too via the Window object.
try {
    $excel.ActiveWindow.ScrollRow = 6
    $excel.ActiveWindow.ScrollColumn = 1
} catch {
}
$ws.Range("G15").Select() | Out-Null

# --- Row 12 ("Total Classes") --------------------------------------------
# The 10th (col F) and 17th (col G) counts are corrected to 0, halving the
# monthly total from 12 to 6, and the %age formula is repointed from the
# (self-referential, always-100) I12 to I9 (blank -> 0%).
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 0
$ws.Range("J12").Formula = "=(I9/`$I`$12)*100"

# --- Column E (attendance on the 10th) for students in rows 14-23 --------
# Marked present (3) except for the three students who were absent (0).
$presentRows = 14, 15, 17, 18, 19, 20, 22
$absentRows  = 16, 21, 23
foreach ($r in $presentRows) {
    $ws.Range("E$r").Value = 3
}
foreach ($r in $absentRows) {
    $ws.Range("E$r").Value = 0
}
